$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top, pushing the existing header row
# (currently row 1) down to row 2 - matching the target location.
$ws.Rows.Item(1).Insert()
